$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manager closing kiosk: update pizza ingredient inventory counts
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = -1231
$ws.Range("C2").Value = -1231
$ws.Range("D2").Value = -1231
$ws.Range("G2").Value = 1000
